$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 205.375
$ws.Range("I12").Value = 205.375
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 205.375
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -35.375
$ws.Range("H18").Value = 2050
$ws.Range("I18").Value = 2050
$ws.Range("K18").Value = 2050
$ws.Range("M18").Value = -1766
$ws.Range("H31").Value = 2662.8572
$ws.Range("I31").Value = 2662.8572
$ws.Range("K31").Value = 7988.571599999999
$ws.Range("M31").Value = -7758.571599999999
$ws.Range("H40").Value = 3999
$ws.Range("I40").Value = 3999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3999
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3824
$ws.Range("H98").Value = 2224.8333
$ws.Range("I98").Value = 1885.1177
$ws.Range("K98").Value = 1885.1177
$ws.Range("M98").Value = -387.1177
$ws.Range("H122").Value = 2224.8333
$ws.Range("I122").Value = 1885.1177
$ws.Range("K122").Value = 5655.3531
$ws.Range("M122").Value = -3205.3531
$ws.Range("H132").Value = 3869.7058
$ws.Range("I132").Value = 3021.25
$ws.Range("K132").Value = 9063.75
$ws.Range("M132").Value = -6533.75
$ws.Range("H138").Value = 4004.1052
$ws.Range("J138").Value = 4420.8887
$ws.Range("L138").Value = 13262.6661
$ws.Range("N138").Value = -23542.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31145.7
$ws.Range("I32").Value = 34419.21
$ws.Range("K32").Value = 34419.21
$ws.Range("M32").Value = -34132.21
$ws.Range("H35").Value = 10509.25
$ws.Range("J35").Value = 16519
$ws.Range("L35").Value = 16519
$ws.Range("N35").Value = -17331
$ws.Range("H37").Value = 6999
$ws.Range("I37").Value = 6999
$ws.Range("K37").Value = 6999
$ws.Range("M37").Value = -6726
$ws.Range("H61").Value = 5443.696
$ws.Range("I61").Value = 5443.696
$ws.Range("K61").Value = 5443.696
$ws.Range("M61").Value = -5231.696
$ws.Range("H74").Value = 50090.145
$ws.Range("I74").Value = 55087.684
$ws.Range("K74").Value = 55087.684
$ws.Range("M74").Value = -54213.684
$ws.Range("H77").Value = 50090.145
$ws.Range("I77").Value = 55087.684
$ws.Range("K77").Value = 275438.42
$ws.Range("M77").Value = -271070.42
$ws.Range("H132").Value = 27049.268
$ws.Range("I132").Value = 28301.795
$ws.Range("K132").Value = 84905.38499999999
$ws.Range("M132").Value = -82375.38499999999
$ws.Range("H136").Value = 5443.696
$ws.Range("I136").Value = 5443.696
$ws.Range("K136").Value = 16331.088
$ws.Range("M136").Value = -13781.088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 39439
$ws.Range("I26").Value = 8888
$ws.Range("K26").Value = 8888
$ws.Range("M26").Value = -8596
$ws.Range("H64").Value = 2199.4211
$ws.Range("I64").Value = 1458.2
$ws.Range("J64").Value = 3023
$ws.Range("K64").Value = 1458.2
$ws.Range("L64").Value = 3023
$ws.Range("M64").Value = -1233.2
$ws.Range("N64").Value = -3473
$ws.Range("H67").Value = 2199.4211
$ws.Range("I67").Value = 1458.2
$ws.Range("J67").Value = 3023
$ws.Range("K67").Value = 1458.2
$ws.Range("L67").Value = 3023
$ws.Range("M67").Value = -678.2
$ws.Range("N67").Value = -4583
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H99").Value = 2212.8823
$ws.Range("I99").Value = 2418.7778
$ws.Range("K99").Value = 2418.7778
$ws.Range("M99").Value = -920.7777999999998
$ws.Range("H107").Value = 2032.9584
$ws.Range("I107").Value = 1361.2941
$ws.Range("K107").Value = 1361.2941
$ws.Range("M107").Value = 558.7058999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 816.8
$ws.Range("I16").Value = 746.6667
$ws.Range("K16").Value = 746.6667
$ws.Range("M16").Value = -459.6667
$ws.Range("H31").Value = 3613.7144
$ws.Range("I31").Value = 2165.3333
$ws.Range("J31").Value = 4700
$ws.Range("K31").Value = 2165.3333
$ws.Range("L31").Value = 4700
$ws.Range("M31").Value = -1870.3333
$ws.Range("N31").Value = -5290
$ws.Range("H34").Value = 3613.7144
$ws.Range("I34").Value = 2165.3333
$ws.Range("J34").Value = 4700
$ws.Range("K34").Value = 2165.3333
$ws.Range("L34").Value = 4700
$ws.Range("M34").Value = -1963.3333
$ws.Range("N34").Value = -5104
$ws.Range("H41").Value = 16750
$ws.Range("I41").Value = 16750
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 16750
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -16322
$ws.Range("H52").Value = 52999.668
$ws.Range("J52").Value = 57999.6
$ws.Range("L52").Value = 57999.6
$ws.Range("N52").Value = -58587.6
$ws.Range("H70").Value = 61997.5
$ws.Range("J70").Value = 73995
$ws.Range("L70").Value = 73995
$ws.Range("N70").Value = -74625
$ws.Range("H73").Value = 61997.5
$ws.Range("J73").Value = 73995
$ws.Range("L73").Value = 73995
$ws.Range("N73").Value = -76179
$ws.Range("H99").Value = 144444.86
$ws.Range("I99").Value = 201217.8
$ws.Range("K99").Value = 201217.8
$ws.Range("M99").Value = -199719.8
$ws.Range("H113").Value = 816.8
$ws.Range("I113").Value = 746.6667
$ws.Range("K113").Value = 746.6667
$ws.Range("M113").Value = 1423.3333
$ws.Range("H122").Value = 1534
$ws.Range("I122").Value = 1382.7778
$ws.Range("J122").Value = 1874.25
$ws.Range("K122").Value = 4148.3334
$ws.Range("L122").Value = 5622.75
$ws.Range("M122").Value = -1698.3334
$ws.Range("N122").Value = -10522.75
$ws.Range("H126").Value = 144444.86
$ws.Range("I126").Value = 201217.8
$ws.Range("K126").Value = 603653.3999999999
$ws.Range("M126").Value = -601183.3999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 499.5
$ws.Range("I35").Value = 499
$ws.Range("K35").Value = 1497
$ws.Range("M35").Value = -1209
$ws.Range("H119").Value = 1910.4
$ws.Range("I119").Value = 1910.4
$ws.Range("K119").Value = 5731.200000000001
$ws.Range("M119").Value = -893.2000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16331.1
$ws.Range("J24").Value = 16331.1
$ws.Range("L24").Value = 16331.1
$ws.Range("N24").Value = -16677.1
$ws.Range("H42").Value = 99290
$ws.Range("J42").Value = 99290
$ws.Range("L42").Value = 99290
$ws.Range("N42").Value = -100260
$ws.Range("H115").Value = 99290
$ws.Range("J115").Value = 99290
$ws.Range("L115").Value = 99290
$ws.Range("N115").Value = -101640
$ws.Range("H132").Value = 114989.2
$ws.Range("I132").Value = 161085
$ws.Range("J132").Value = 7432.3335
$ws.Range("K132").Value = 483255
$ws.Range("L132").Value = 22297.0005
$ws.Range("M132").Value = -480725
$ws.Range("N132").Value = -27357.0005
$ws.Range("H136").Value = 23299.555
$ws.Range("J136").Value = 23299.555
$ws.Range("L136").Value = 69898.66500000001
$ws.Range("N136").Value = -74998.66500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 76585.13
$ws.Range("I22").Value = 370629.66
$ws.Range("J22").Value = 3074
$ws.Range("K22").Value = 370629.66
$ws.Range("L22").Value = 3074
$ws.Range("M22").Value = -370334.66
$ws.Range("N22").Value = -3664
$ws.Range("H27").Value = 76585.13
$ws.Range("I27").Value = 370629.66
$ws.Range("J27").Value = 3074
$ws.Range("K27").Value = 370629.66
$ws.Range("L27").Value = 3074
$ws.Range("M27").Value = -370522.66
$ws.Range("N27").Value = -3288
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H61").Value = 1533.8422
$ws.Range("I61").Value = 935.6
$ws.Range("K61").Value = 935.6
$ws.Range("M61").Value = -733.6
$ws.Range("H113").Value = 1533.8422
$ws.Range("I113").Value = 935.6
$ws.Range("K113").Value = 935.6
$ws.Range("M113").Value = 1234.4
$ws.Range("H122").Value = 4280.9375
$ws.Range("I122").Value = 3358.5715
$ws.Range("J122").Value = 4998.3335
$ws.Range("K122").Value = 10075.7145
$ws.Range("L122").Value = 14995.0005
$ws.Range("M122").Value = -7625.7145
$ws.Range("N122").Value = -19895.0005
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 27155.18
$ws.Range("I132").Value = 31790.342
$ws.Range("J132").Value = 6039.4443
$ws.Range("K132").Value = 95371.026
$ws.Range("L132").Value = 18118.3329
$ws.Range("M132").Value = -92841.026
$ws.Range("N132").Value = -23178.3329
$ws.Range("H136").Value = 4338.567
$ws.Range("I136").Value = 2781.1667
$ws.Range("J136").Value = 6674.6665
$ws.Range("K136").Value = 8343.500100000001
$ws.Range("L136").Value = 20023.9995
$ws.Range("M136").Value = -5793.500100000001
$ws.Range("N136").Value = -25123.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 16668634
$ws.Range("I96").Value = 2360.4
$ws.Range("K96").Value = 2360.4
$ws.Range("M96").Value = -987.4000000000001
$ws.Range("H123").Value = 99997
$ws.Range("J123").Value = 99997
$ws.Range("L123").Value = 99997
$ws.Range("N123").Value = -109797
$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -100119
